# This script regenerates the "K" column (column G) of the save_data
# sheet with freshly computed values (commit: "regen save_data to use K
# instead of Strike#, regen std/mean, calc and write s_vals").
#
# Only column G ("K") values change; every other column is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of worksheet row number -> new value for column G ("K").
$updates = @(
    @{ Row = 2;  Value = 0 },
    @{ Row = 3;  Value = 1 },
    @{ Row = 4;  Value = 2 },
    @{ Row = 5;  Value = 3 },
    @{ Row = 6;  Value = 1 },
    @{ Row = 7;  Value = 1 },
    @{ Row = 8;  Value = 2 },
    @{ Row = 9;  Value = 1 },
    @{ Row = 10; Value = 0 },
    @{ Row = 11; Value = 1 },
    @{ Row = 12; Value = 3 },
    @{ Row = 13; Value = 3 },
    @{ Row = 14; Value = 1 },
    @{ Row = 15; Value = 2 },
    @{ Row = 16; Value = 0 },
    @{ Row = 17; Value = 1 },
    @{ Row = 18; Value = 2 },
    @{ Row = 19; Value = 1 },
    @{ Row = 20; Value = 1 },
    @{ Row = 21; Value = 0 },
    @{ Row = 22; Value = 0 },
    @{ Row = 23; Value = 3 },
    @{ Row = 24; Value = 1 },
    @{ Row = 25; Value = 3 },
    @{ Row = 26; Value = 0 },
    @{ Row = 27; Value = 2 },
    @{ Row = 28; Value = 1 },
    @{ Row = 29; Value = 1 },
    @{ Row = 30; Value = 0 },
    @{ Row = 31; Value = 1 },
    @{ Row = 33; Value = 3 },
    @{ Row = 34; Value = 1 },
    @{ Row = 35; Value = 2 },
    @{ Row = 36; Value = 1 },
    @{ Row = 37; Value = 0 },
    @{ Row = 38; Value = 1 },
    @{ Row = 39; Value = 0 },
    @{ Row = 40; Value = 0 },
    @{ Row = 41; Value = 2 },
    @{ Row = 42; Value = 3 },
    @{ Row = 43; Value = 1 },
    @{ Row = 44; Value = 2 },
    @{ Row = 45; Value = 1 },
    @{ Row = 46; Value = 0 },
    @{ Row = 47; Value = 2 },
    @{ Row = 48; Value = 1 },
    @{ Row = 49; Value = 1 },
    @{ Row = 50; Value = 1 },
    @{ Row = 51; Value = 2 },
    @{ Row = 52; Value = 3 },
    @{ Row = 53; Value = 1 },
    @{ Row = 54; Value = 1 },
    @{ Row = 55; Value = 1 },
    @{ Row = 57; Value = 1 },
    @{ Row = 58; Value = 2 },
    @{ Row = 59; Value = 0 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Value
}
